$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Re-use the pre-existing "hol_date" number format (mm-dd-yy / numFmtId 14,
# already style index 1 on the original column C) for the new hol_date
# column (B) BEFORE anything else is touched, so the engine keeps reusing
# the same cellXf instead of minting a duplicate.
# ---------------------------------------------------------------------
$ws.Range("C2").Copy() | Out-Null
$ws.Range("B2:B10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# The old hol_date column (C) becomes hol_name, which is plain text --
# clear its leftover date formatting.
$ws.Range("C2:C10").Clear() | Out-Null

# New cr_dtimes column (H) uses a single "mm:ss.0" (numFmtId 47) format,
# shared by every row.
$ws.Range("H2").NumberFormat = "mm:ss.0"
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H2:H10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---- Header row -----------------------------------------------------
$headers = @("regcntr_id","hol_date","hol_name","hol_reason","lang_code","is_active","cr_by","cr_dtimes","upd_by","upd_dtimes","is_deleted","del_dtimes")
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# ---- Data rows --------------------------------------------------------
# date serial | hol_name (= hol_reason)
$rows = @(
    @(44927, "Jour de l'an"),
    @(45037, "FÃªte de Ramadan"),
    @(45047, "FÃªte du travail"),
    @(45071, "Anniversaire de l'Union Africaine "),
    @(45107, "Tabaski"),
    @(45153, "Assomption"),
    @(45196, "Maouloud"),
    @(45201, "FÃªte de l'indÃ©pendance"),
    @(45285, "Jour de Noel")
)

$rowNum = 2
foreach ($item in $rows) {
    $date = $item[0]
    $name = $item[1]

    $ws.Cells.Item($rowNum, 1).Value = 10001
    $ws.Cells.Item($rowNum, 2).Value = $date
    $ws.Cells.Item($rowNum, 3).Value = $name
    $ws.Cells.Item($rowNum, 4).Value = $name
    $ws.Cells.Item($rowNum, 5).Value = "fra"
    $ws.Cells.Item($rowNum, 6).Value = $true
    $ws.Cells.Item($rowNum, 7).Value = "superadmin"
    $ws.Cells.Item($rowNum, 8).Value = 45079.578724606479
    $ws.Cells.Item($rowNum, 9).Value = "NULL"
    $ws.Cells.Item($rowNum, 10).Value = "NULL"
    $ws.Cells.Item($rowNum, 11).Value = $false
    $ws.Cells.Item($rowNum, 12).Value = "NULL"

    $rowNum++
}

# ---- Selection ------------------------------------------------------------
$ws.Range("D17").Select() | Out-Null
